$wb = $excel.ActiveWorkbook

# Add the new "forms" worksheet after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "forms"

# Header row
$ws.Range("A1").Value = "FormNameText"
$ws.Range("B1").Value = "IntroText"
$ws.Range("C1").Value = "CompleteText"

# Data rows
$ws.Range("A2").Value = "TOM"
$ws.Range("B2").Value = "Testing123"
$ws.Range("C2").Value = "GSDFSDCSXCC"

$ws.Range("A3").Value = "David"
$ws.Range("B3").Value = "Testing345"
$ws.Range("C3").Value = "ASDASDASD"

$ws.Range("A4").Value = "Mukta"
$ws.Range("B4").Value = "Testing567"
$ws.Range("C4").Value = "GDFSDFFFF"

# Style the Intro/Complete header cells with a distinct font.
$ws.Range("B1:C1").Font.Name = "Consolas"
$ws.Range("B1:C1").Font.Size = 10
$ws.Range("B1:C1").Font.Color = 4079210

# Column widths to fit the new content (best-fit-style autosize).
$ws.Columns.Item(1).ColumnWidth = 12.917
$ws.Columns.Item(2).ColumnWidth = 9.084
$ws.Columns.Item(3).ColumnWidth = 12.25

# Move the selection like the source workbook (C6, below the data).
$ws.Range("C6").Select()

# "contacts" is no longer the active tab; "forms" is.
$wb.Worksheets.Item("forms").Activate()
